# corrected file endings in Excel files
# Change file extension in column A (rows 2-7) from .xmi to .tsv

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "rwk1_digbib_300-1.tsv"
$ws.Range("A3").Value = "rwk1_digbib_1039-1.tsv"
$ws.Range("A4").Value = "rwk1_digbib_1057-1.tsv"
$ws.Range("A5").Value = "rwk1_mkhz_2778-1.tsv"
$ws.Range("A6").Value = "rwk1_mkhz_6147-1.tsv"
$ws.Range("A7").Value = "rwk1_mkhz_6263-1.tsv"

# Update selection to match the committed state (A1:A1048576 selected, whole column A)
$ws.Columns("A").Select()

# Set column A width to fit the new filenames (target stored width ~21.7109375, i.e. ~21 characters wide)
$ws.Columns("A").ColumnWidth = 20.9
